$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.614.53'
$ws.Range('D3').Value = '2.492.80'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.69'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.74'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').Value = '  -1.47%  '
$ws.Range('D9').Value = '2.493.46'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.140'
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.09'
$ws.Range('E12').Value = '  -1.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.341'
$ws.Range('E13').Value = '  -2.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.27'
$ws.Range('E14').Value = '  -3.04%  '
$ws.Range('D15').Value = '2.948.03'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').Value = '67.499.59'
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '2.487.30'
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.82'
$ws.Range('E19').Value = '  +3.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.00'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '366.60'
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.12'
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.55'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.33'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  -5.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.97'
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = '2.623.02'
$ws.Range('E29').Value = '  -2.18%  '
$ws.Range('D30').Value = '0.0₃0959'
$ws.Range('E30').Value = '  -3.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.31'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '531.53'
$ws.Range('E32').Value = '  -1.42%  '
$ws.Range('E33').Value = '  -4.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.86'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.77'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('E38').Value = '  -3.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.68'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.64'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('E43').Value = '  -1.53%  '
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.50'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '144.89'
$ws.Range('E46').Value = '  -3.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.67'
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.546'
$ws.Range('E48').Value = '  -3.27%  '
$ws.Range('D49').Value = '0.0₆0273'
$ws.Range('E49').Value = '  -2.53%  '
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('E51').Value = '  -2.00%  '
